# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 91 ("Feria Lagunitas de
# Puerto Montt" / "Sandia"), pushing the previously existing rows 91-123
# down to 92-124 (dimension grows from A1:R123 to A1:R124).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 91; this shifts rows 91..123
# down to 92..124 and keeps all of their data/styles intact.
$ws.Rows(91).Insert()

# Populate the newly inserted row 91 with the new record.
$ws.Range("A91").Value = 4
$ws.Range("B91").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C91").Value = "Los Lagos"
$ws.Range("D91").Value = 44488
$ws.Range("E91").Value = 10
$ws.Range("F91").Value = 100112028
$ws.Range("G91").Value = "Sandia"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 600
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 1000
$ws.Range("M91").Value = 1000
$ws.Range("N91").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O91").Value = "Perú"
$ws.Range("P91").Value = 1000
$ws.Range("Q91").Value = 1
$ws.Range("R91").Value = "Hortaliza"
